$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("dSF") updates per row, based on repulled/pushed data.
$updates = @{
    2  = -2
    3  = -2
    5  = 4
    6  = 3
    7  = -2
    8  = -2
    9  = -5
    10 = 7
    12 = -1
    14 = 6
    15 = -1
    16 = 3
    17 = 2
    18 = -1
    19 = -2
    20 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
